# modify the detailUserActivity for the user info
# - Populate Sheet2 with the SecretManager / DetailUserActivity note content
# - Make Sheet2 the active (selected) sheet/tab instead of "db"
# - Update the remembered selection on "db" (B17 -> B18) and on Sheet2 (-> B7)

$wb  = $excel.ActiveWorkbook
$db  = $wb.Worksheets.Item("db")
$ws2 = $wb.Worksheets.Item("Sheet2")

# New content for Sheet2
$ws2.Range("B2").Value = "SecretManager와 DetailUserActivity간의 통신"
$ws2.Range("B3").Value = "putExtra 사용"
$ws2.Range("B4").Value = "name"
$ws2.Range("B5").Value = "tel"
$ws2.Range("B7").Value = "단순하게 이름과 전화번호만 필요하다."

# "db" keeps its old selection, just moved down one row (B17 -> B18)
$db.Range("B18").Select() | Out-Null

# Sheet2 becomes the active tab, with its own new selection (B7)
$ws2.Activate() | Out-Null
$ws2.Range("B7").Select() | Out-Null
